$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.123.83"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.669.17"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'210.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'0.5197"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D8").Value = "'0.2611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "'0.06327"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'21.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "'0.07544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.673.51"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "'4.416"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "'0.5430"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "'0.000008033"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'66.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "26.178.10"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'4.736"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "'187.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "'10.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").Value = "'6.242"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'149.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'0.1234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").Value = "'7.477"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "'15.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'0.06315"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "'1.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'3.418"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("D33").Value = "'1.648"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "'1.002"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'0.6003"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'2.399"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "'2.765"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").Value = "1.113.73"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("D39").Value = "'0.01612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'6.047"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "'0.8643"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'100.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").Value = "1.822.18"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'55.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").Value = "'8.050"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'0.4238"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'5.900"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
